# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the handback
# files for d3076c5a-7e2e-4bf4-89cf-318aa9b55203 and
# df845ac9-e96b-4090-af84-cde44e5b3f39 are out of date vs. the latest
# source, so both are now "Ready for handoff" instead of
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$errMsgD3076 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0531b1988f7be651ee1d0ac82f1576737f62e40c/e2e/d3076c5a-7e2e-4bf4-89cf-318aa9b55203.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2f3e4b9e7a622d44edb6794d567f2428c85bdbc/e2e/d3076c5a-7e2e-4bf4-89cf-318aa9b55203.md."
$errMsgDf845 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0531b1988f7be651ee1d0ac82f1576737f62e40c/e2e/df845ac9-e96b-4090-af84-cde44e5b3f39.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2f3e4b9e7a622d44edb6794d567f2428c85bdbc/e2e/df845ac9-e96b-4090-af84-cde44e5b3f39.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-15 22:25:20"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-15 22:25:20"

# ---- zh-cn sheet ----
# Note: Excel's ColumnWidth setter pads the stored OOXML width by 5/6 of a
# character (the default-font column-width quirk), so ask for 40 - 5/6 to
# land on a stored width of exactly 40.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("H4").Value = "2016-08-15 22:25:16"
$wsZhCn.Range("P4").Value = $errMsgD3076
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("H5").Value = "2016-08-15 22:25:16"
$wsZhCn.Range("P5").Value = $errMsgDf845
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("H4").Value = "2016-08-15 22:25:20"
$wsDeDe.Range("P4").Value = $errMsgD3076
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("H5").Value = "2016-08-15 22:25:20"
$wsDeDe.Range("P5").Value = $errMsgDf845
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
